$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly-annotated columns E:J for rows 2-12 ---
$data = @{
    2  = @(2, 2, 1, 2, 2, 2)
    3  = @(2, 2, 1, 1, 2, 2)
    4  = @(2, 2, 1, 1, 1, 2)
    5  = @(2, 2, 1, 2, 1, 2)
    6  = @(2, 2, 1, 1, 2, 2)
    7  = @(2, 2, 2, 2, 2, 2)
    8  = @(2, 2, 1, 2, 1, 2)
    9  = @(2, 2, 1, 2, 2, 2)
    10 = @(2, 2, 1, 1, 2, 2)
    11 = @(2, 2, 1, 1, 2, 2)
    12 = @(2, 2, 1, 2, 1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # E=5 .. J=10
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# --- View changes: zoom, freeze header row, scroll, final selection ---
$excel.ActiveWindow.Zoom = 85

$ws.Rows("2:2").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7

$ws.Range("E13").Select()
